$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 'last updated' timestamp in A1
$ws.Range('A1').Value = 'Datos actualizados a 7 de Abril de 2020 a las 17:22'

# Refresh per-country COVID data (values updated + table re-sorted by Casos totales desc)
$data = @(
    @('Estados Unidos', 377322, 10318, 19877, 345672, 9015, 902, 11773),
    @('España', 140510, 3835, 43208, 83504, 7069, 457, 13798),
    @('Italia', 132547, 0, 22837, 93187, 3898, 0, 16523),
    @('Alemania', 104592, 1217, 36081, 66657, 4895, 44, 1854),
    @('Francia', 98010, 0, 17250, 71849, 7072, 0, 8911),
    @('China', 81740, 32, 77167, 1242, 211, 0, 3331),
    @('Iran', 62589, 2089, 27039, 31678, 3987, 133, 3872),
    @('Reino Unido', 55242, 3634, 135, 48948, 1559, 786, 6159),
    @('Turquia', 30217, 0, 1326, 28242, 1415, 0, 649),
    @('Suiza', 22242, 585, 8056, 13375, 391, 46, 811),
    @('Belgica', 22194, 1380, 4157, 16002, 1260, 403, 2035),
    @('Paises Bajos', 19580, 777, 250, 17229, 1424, 234, 2101),
    @('Canada', 17046, 379, 3794, 12907, 426, 22, 345),
    @('Austria', 12592, 295, 4046, 8303, 243, 23, 243),
    @('Portugal', 12442, 712, 184, 11913, 271, 34, 345),
    @('Brasil', 12345, 162, 127, 11637, 296, 17, 581),
    @('Corea del Sur', 10331, 47, 6694, 3445, 55, 6, 192),
    @('Israel', 9006, 102, 683, 8263, 153, 3, 60),
    @('Suecia', 7693, 487, 205, 6897, 640, 114, 591),
    @('Rusia', 7497, 1154, 494, 6945, 8, 11, 58),
    @('Australia', 5919, 24, 2547, 3324, 93, 3, 48),
    @('Noruega', 5866, 1, 32, 5746, 83, 12, 88),
    @('Irlanda', 5364, 0, 25, 5165, 165, 0, 174),
    @('Chile', 5116, 301, 898, 4175, 337, 6, 43),
    @('Dinamarca', 5071, 390, 1491, 3377, 127, 16, 203),
    @('India', 4911, 133, 382, 4392, 0, 1, 137),
    @('Chequia', 4828, 6, 127, 4621, 86, 2, 80),
    @('Polonia', 4532, 119, 191, 4230, 50, 4, 111),
    @('Rumania', 4417, 360, 460, 3760, 274, 21, 197),
    @('Pakistan', 4005, 239, 429, 3521, 28, 2, 55),
    @('Malasia', 3963, 170, 1321, 2579, 92, 1, 63),
    @('Japon', 3906, 0, 592, 3222, 79, 0, 92),
    @('Filipinas', 3764, 104, 84, 3503, 1, 14, 177),
    @('Ecuador', 3747, 0, 100, 3456, 156, 0, 191),
    @('Luxemburgo', 2843, 0, 500, 2302, 33, 0, 41),
    @('Arabia Saudita', 2795, 190, 615, 2139, 41, 3, 41),
    @('Indonesia', 2738, 247, 204, 2313, 0, 12, 221),
    @('Peru', 2561, 0, 997, 1472, 89, 0, 92),
    @('Serbia', 2447, 247, 118, 2268, 109, 3, 61),
    @('Mexico', 2439, 296, 633, 1681, 89, 31, 125),
    @('Finlandia', 2308, 132, 300, 1974, 81, 7, 34),
    @('Tailandia', 2258, 38, 824, 1407, 30, 1, 27),
    @('Panama', 2100, 0, 14, 2031, 88, 0, 55),
    @('Emiratos Arabes Unidos', 2076, 0, 167, 1898, 1, 0, 11),
    @('Republica Dominicana', 1956, 128, 33, 1825, 147, 12, 98),
    @('Catar', 1832, 0, 131, 1697, 37, 0, 4),
    @('Grecia', 1832, 77, 269, 1482, 90, 2, 81),
    @('Sudafrica', 1749, 63, 95, 1641, 7, 1, 13),
    @('Argentina', 1628, 0, 338, 1235, 96, 2, 55),
    @('Islandia', 1586, 24, 559, 1021, 11, 0, 6),
    @('Colombia', 1579, 0, 88, 1445, 50, 0, 46),
    @('Ucrania', 1462, 143, 28, 1389, 16, 7, 45),
    @('Argelia', 1423, 0, 90, 1160, 46, 0, 173),
    @('Singapur', 1375, 0, 344, 1025, 25, 0, 6),
    @('Egipto', 1322, 0, 259, 978, 0, 0, 85),
    @('Croacia', 1282, 60, 167, 1097, 35, 2, 18),
    @('Nueva Zelanda', 1160, 54, 241, 918, 14, 0, 1),
    @('Estonia', 1149, 41, 69, 1059, 12, 2, 21),
    @('Marruecos', 1141, 21, 88, 970, 1, 3, 83),
    @('Eslovenia', 1059, 38, 102, 921, 30, 6, 36),
    @('Moldavia', 1056, 91, 40, 994, 80, 3, 22),
    @('Irak', 1031, 0, 344, 623, 0, 0, 64),
    @('Hong Kong', 936, 21, 236, 696, 12, 0, 4),
    @('Lituania', 880, 37, 8, 857, 11, 0, 15),
    @('Bielorrusia', 861, 161, 54, 794, 31, 0, 13),
    @('Armenia', 853, 20, 87, 758, 30, 0, 8),
    @('Hungria', 817, 73, 71, 699, 17, 9, 47),
    @('Barein', 811, 55, 458, 349, 4, 0, 4),
    @('Bosnia y Herzegovina', 754, 80, 68, 653, 4, 4, 33),
    @('Kuwait', 743, 78, 105, 637, 23, 0, 1),
    @('Azerbaiyan', 717, 76, 44, 665, 23, 1, 8),
    @('Crucero', 712, 0, 619, 82, 10, 0, 11),
    @('Kazajistan', 685, 23, 50, 629, 16, 0, 6),
    @('Camerun', 658, 0, 17, 632, 0, 0, 9),
    @('Republica de Macedonia', 599, 29, 30, 543, 15, 3, 26),
    @('Tunez', 596, 0, 25, 549, 39, 0, 22),
    @('Eslovaquia', 581, 47, 13, 566, 3, 0, 2),
    @('Bulgaria', 565, 16, 42, 501, 26, 0, 22),
    @('Letonia', 548, 6, 16, 530, 5, 1, 2),
    @('Libano', 548, 7, 60, 469, 27, 0, 19),
    @('Principado de Andorra', 525, 0, 31, 473, 12, 0, 21),
    @('Uzbekistan', 504, 47, 30, 472, 8, 0, 2),
    @('Republica de Chipre', 494, 29, 45, 440, 11, 0, 9),
    @('Costa Rica', 467, 0, 18, 447, 14, 0, 2),
    @('Afganistan', 423, 56, 18, 391, 0, 3, 14),
    @('Uruguay', 415, 0, 123, 286, 14, 0, 6),
    @('Albania', 383, 6, 131, 230, 7, 1, 22),
    @('Taiwan', 376, 3, 61, 310, 0, 0, 5),
    @('Oman', 371, 40, 67, 302, 3, 0, 2),
    @('Burkina Faso', 364, 0, 108, 238, 0, 0, 18),
    @('Cuba', 363, 13, 18, 336, 12, 0, 9),
    @('Reunion', 349, 0, 40, 309, 4, 0, 0),
    @('Jordania', 349, 0, 126, 217, 5, 0, 6),
    @('Costa de Marfil', 323, 0, 41, 279, 0, 0, 3),
    @('Honduras', 305, 7, 6, 277, 10, 0, 22),
    @('Malta', 293, 52, 5, 288, 4, 0, 0),
    @('Ghana', 287, 73, 31, 251, 2, 0, 5),
    @('San Marino', 277, 0, 35, 210, 14, 0, 32),
    @('Mauricio', 268, 24, 8, 253, 3, 0, 7),
    @('Estado de Palestina', 261, 7, 42, 218, 0, 0, 1),
    @('Niger', 253, 0, 26, 217, 0, 0, 10),
    @('Vietnam', 249, 4, 123, 126, 8, 0, 0),
    @('Montenegro', 241, 8, 4, 235, 7, 0, 2),
    @('Nigeria', 238, 0, 35, 198, 2, 0, 5),
    @('Senegal', 237, 11, 105, 130, 1, 0, 2),
    @('Kirguistan', 228, 12, 33, 191, 5, 0, 4),
    @('Georgia', 195, 7, 45, 147, 6, 1, 3),
    @('Bolivia', 194, 11, 2, 178, 3, 3, 14),
    @('Sri Lanka', 185, 7, 42, 137, 5, 1, 6),
    @('Islas Feroe', 184, 1, 120, 64, 2, 0, 0),
    @('Kenia', 172, 14, 7, 159, 2, 0, 6),
    @('Venezuela', 165, 0, 65, 93, 6, 0, 7),
    @('Mayotte', 164, 0, 15, 147, 3, 0, 2),
    @('Banglades', 164, 41, 33, 114, 1, 5, 17),
    @('Consejo Danes para los Refugiados', 161, 0, 5, 138, 0, 0, 18),
    @('Martinica', 151, 0, 50, 97, 20, 0, 4),
    @('Isla de Man', 150, 11, 73, 76, 6, 0, 1),
    @('Guadalupe', 139, 0, 31, 101, 14, 0, 7),
    @('Brunei', 135, 0, 85, 49, 3, 0, 1),
    @('Guinea', 128, 0, 5, 123, 0, 0, 0),
    @('Paraguay', 115, 2, 15, 95, 1, 0, 5),
    @('Camboya', 115, 1, 58, 57, 1, 0, 0),
    @('Gibraltar', 113, 4, 60, 53, 0, 0, 0),
    @('Trinidad yTobago', 106, 1, 1, 97, 0, 0, 8),
    @('Ruanda', 105, 0, 4, 101, 0, 0, 0),
    @('Republica de Yibuti', 90, 0, 9, 81, 0, 0, 0),
    @('Madagascar', 82, 0, 2, 80, 6, 0, 0),
    @('El Salvador', 78, 9, 5, 69, 4, 0, 4),
    @('Monaco', 77, 0, 4, 72, 4, 0, 1),
    @('Guatemala', 77, 7, 17, 57, 3, 0, 3),
    @('Liechtenstein', 77, 0, 55, 21, 0, 0, 1),
    @('Guayana Francesa', 72, 0, 34, 38, 1, 0, 0),
    @('Aruba', 71, 0, 2, 69, 0, 0, 0),
    @('Barbados', 60, 0, 6, 52, 4, 0, 2),
    @('Jamaica', 59, 0, 8, 48, 0, 0, 3),
    @('Togo', 58, 0, 23, 32, 0, 0, 3),
    @('Mali', 56, 9, 12, 39, 0, 0, 5),
    @('Uganda', 52, 0, 0, 52, 0, 0, 0),
    @('Etiopia', 52, 8, 4, 46, 1, 0, 2),
    @('Congo', 45, 0, 2, 38, 0, 0, 5),
    @('Macao', 44, 0, 10, 34, 1, 0, 0),
    @('Polinesia Francesa', 42, 0, 0, 42, 0, 0, 0),
    @('Islas Caimanes', 39, 0, 1, 37, 0, 0, 1),
    @('Puerto Rico', 39, 0, 1, 36, 0, 0, 2),
    @('Zambia', 39, 0, 7, 31, 0, 0, 1),
    @('Bermudas', 39, 0, 17, 20, 0, 0, 2),
    @('San Martin (Parte Holandesa)', 37, 0, 1, 30, 0, 0, 6),
    @('Guinea-Bisau', 33, 15, 0, 33, 0, 0, 0),
    @('Bahamas', 33, 0, 5, 23, 1, 0, 5),
    @('Guam', 32, 0, 0, 31, 0, 0, 1),
    @('San Martin (Parte Francesa)', 32, 0, 7, 23, 6, 0, 2),
    @('Eritrea', 31, 0, 0, 31, 0, 0, 0),
    @('Guyana', 31, 0, 8, 18, 8, 1, 5),
    @('Gabon', 30, 6, 1, 28, 0, 0, 1),
    @('Benin', 26, 0, 5, 20, 0, 0, 1),
    @('Haiti', 25, 1, 0, 24, 0, 0, 1),
    @('Tanzania', 24, 0, 5, 18, 0, 0, 1),
    @('Birmania', 22, 0, 0, 21, 0, 0, 1),
    @('Libia', 19, 0, 1, 17, 0, 0, 1),
    @('Siria', 19, 0, 2, 15, 0, 0, 2),
    @('Maldivas', 19, 0, 13, 6, 0, 0, 0),
    @('Nueva Caledonia', 18, 0, 1, 17, 0, 0, 0),
    @('Islas Virgenes de los Estados Unidos', 17, 0, 0, 17, 0, 0, 0),
    @('Guinea Ecuatorial', 16, 0, 3, 13, 0, 0, 0),
    @('Namibia', 16, 0, 3, 13, 0, 0, 0),
    @('Angola', 16, 0, 2, 12, 0, 0, 2),
    @('Antigua y Barbuda', 15, 0, 0, 15, 1, 0, 0),
    @('Fiyi', 15, 1, 0, 15, 0, 0, 0),
    @('Dominica', 15, 0, 1, 14, 0, 0, 0),
    @('Mongolia', 15, 0, 4, 11, 0, 0, 0),
    @('Laos', 14, 2, 0, 14, 0, 0, 0),
    @('Santa Lucia', 14, 0, 1, 13, 0, 0, 0),
    @('Sudan', 14, 2, 2, 10, 0, 0, 2),
    @('Liberia', 14, 0, 3, 8, 0, 0, 3),
    @('Curazao', 13, 0, 5, 7, 0, 0, 1),
    @('Granada', 12, 0, 0, 12, 2, 0, 0),
    @('San Cristobal y Nieves', 11, 1, 0, 11, 0, 0, 0),
    @('Seychelles', 11, 0, 0, 11, 0, 0, 0),
    @('Groenlandia', 11, 0, 4, 7, 0, 0, 0),
    @('Zimbabue', 10, 0, 0, 9, 0, 0, 1),
    @('Surinam', 10, 0, 0, 9, 0, 0, 1),
    @('Mozambique', 10, 0, 1, 9, 0, 0, 0),
    @('Republica del Chad', 10, 1, 2, 8, 0, 0, 0),
    @('Suazilandia', 10, 0, 4, 6, 0, 0, 0),
    @('Nepal', 9, 0, 1, 8, 0, 0, 0),
    @('Montserrat', 9, 0, 0, 7, 0, 0, 2),
    @('Republica de Africa Central', 8, 0, 0, 8, 0, 0, 0),
    @('Islas Turcas y Caicos', 8, 0, 0, 7, 0, 0, 1),
    @('Malaui', 8, 3, 0, 7, 1, 1, 1),
    @('Santa Sede', 7, 0, 0, 7, 0, 0, 0),
    @('Belice', 7, 0, 0, 6, 1, 0, 1),
    @('San Vicente y las Granadinas', 7, 0, 1, 6, 0, 0, 0),
    @('Somalia', 7, 0, 1, 6, 0, 0, 0),
    @('Cabo Verde', 7, 0, 1, 5, 0, 0, 1),
    @('Sierra Leona', 6, 0, 0, 6, 0, 0, 0),
    @('Botsuana', 6, 0, 0, 5, 0, 0, 1),
    @('Nicaragua', 6, 0, 0, 5, 0, 0, 1),
    @('San Bartolome', 6, 0, 1, 5, 0, 0, 0),
    @('Mauritania', 6, 0, 2, 3, 0, 0, 1),
    @('Butan', 5, 0, 2, 3, 0, 0, 0),
    @('Sahara Occidental', 4, 0, 0, 4, 0, 0, 0),
    @('Santo Tome y Principe', 4, 0, 0, 4, 0, 0, 0),
    @('Gambia', 4, 0, 2, 1, 0, 0, 1),
    @('Islas Virgenes Britanicas', 3, 0, 0, 3, 0, 0, 0),
    @('Anguila', 3, 0, 0, 3, 0, 0, 0),
    @('Burundi', 3, 0, 0, 3, 0, 0, 0),
    @('Papua Nueva Guinea', 2, 0, 0, 2, 0, 0, 0),
    @('Bonaire, San Eustaquio y Saba', 2, 0, 0, 2, 0, 0, 0),
    @('Islas Malvinas', 2, 0, 0, 2, 0, 0, 0),
    @('San Pedro y Miquelon', 1, 0, 0, 1, 0, 0, 0),
    @('Sudan del Sur', 1, 0, 0, 1, 0, 0, 0),
    @('Timor Oriental', 1, 0, 0, 1, 0, 0, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 4 + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
}

Write-Host "Updated $($data.Length) country rows"
